$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the nodes that have not been created yet:
#  - D1 : "(Sibiu, 140, 50, 11)"  (Sibiu node not yet created)
#  - C2 : "(Lugo, 111, 60, 9)"    (Lugo node not yet created)
$ws.Range("D1").ClearContents()
$ws.Range("C2").ClearContents()
